# Fix typo in the shared string used by Sheet2!A1 (missing decimal point
# in the regression equation: "01411" -> "0.1411").
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = '"y(cm/s) = 0.1411x(RPM) -1.7937"'

# Sheet2 was previously the active/selected tab with F9 selected; flip the
# active tab back to Sheet1 and leave Sheet2's remembered selection at B4.
$ws2.Activate()
$ws2.Range("B4").Select()

$ws1.Activate()
$ws1.Range("D12").Select()
